$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 133
$ws.Range("A133").Value = 131
$ws.Range("B133").Value = 6553761
$ws.Range("C133").Value = 'Serbia Super Liga'
$ws.Range("D133").Value = 'Serbia Super Liga'
$ws.Range("E133").Value = 45068.54166666666
$ws.Range("F133").Value = 'FK Radnik Surdulica'
$ws.Range("G133").Value = 'FK Kolubara'
$ws.Range("H133").Value = 2
$ws.Range("I133").Value = 1
$ws.Range("J133").Value = 'H'
$ws.Range("K133").Value = 1.333
$ws.Range("L133").Value = 4.5
$ws.Range("M133").Value = 7
$ws.Range("N133").Value = 1.4
$ws.Range("O133").Value = 4.333
$ws.Range("P133").Value = 6.5
$ws.Range("Q133").Value = -1.25
$ws.Range("R133").Value = 1.975
$ws.Range("S133").Value = 1.825
$ws.Range("T133").Value = 2.5
$ws.Range("U133").Value = 1.875
$ws.Range("V133").Value = 1.925
$ws.Range("W133").Value = 0.3999999999999999
$ws.Range("X133").Value = -1
$ws.Range("Y133").Value = -1
$ws.Range("Z133").Value = -0.5
$ws.Range("AA133").Value = 0.4125
$ws.Range("AB133").Value = 0.875
$ws.Range("AC133").Value = -1

# Row 134
$ws.Range("A134").Value = 132
$ws.Range("B134").Value = 6553338
$ws.Range("C134").Value = 'Serbia Super Liga'
$ws.Range("D134").Value = 'Serbia Super Liga'
$ws.Range("E134").Value = 45068.54166666666
$ws.Range("F134").Value = 'Radnicki Nis'
$ws.Range("G134").Value = 'Spartak Subotica'
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 1
$ws.Range("J134").Value = 'A'
$ws.Range("K134").Value = 1.85
$ws.Range("L134").Value = 3.2
$ws.Range("M134").Value = 3.8
$ws.Range("N134").Value = 2.25
$ws.Range("O134").Value = 3.1
$ws.Range("P134").Value = 3
$ws.Range("Q134").Value = -0.25
$ws.Range("R134").Value = 1.975
$ws.Range("S134").Value = 1.825
$ws.Range("T134").Value = 2.25
$ws.Range("U134").Value = 1.825
$ws.Range("V134").Value = 1.975
$ws.Range("W134").Value = -1
$ws.Range("X134").Value = -1
$ws.Range("Y134").Value = 2
$ws.Range("Z134").Value = -1
$ws.Range("AA134").Value = 0.825
$ws.Range("AB134").Value = -1
$ws.Range("AC134").Value = 0.9750000000000001

# Row 161
$ws.Range("A161").Value = 159
$ws.Range("B161").Value = 7032917
$ws.Range("C161").Value = 'Serbia Super Liga'
$ws.Range("D161").Value = 'Serbia Super Liga'
$ws.Range("E161").Value = 45150.57986111111
$ws.Range("F161").Value = 'FK Backa Topola'
$ws.Range("G161").Value = 'FK Radnicki 1923'
$ws.Range("H161").Value = 1
$ws.Range("I161").Value = 0
$ws.Range("J161").Value = 'H'
$ws.Range("K161").Value = 1.5
$ws.Range("L161").Value = 3.75
$ws.Range("M161").Value = 6.5
$ws.Range("N161").Value = 1.444
$ws.Range("O161").Value = 4
$ws.Range("P161").Value = 6.5
$ws.Range("Q161").Value = -1.25
$ws.Range("R161").Value = 1.975
$ws.Range("S161").Value = 1.825
$ws.Range("T161").Value = 2.75
$ws.Range("U161").Value = 1.95
$ws.Range("V161").Value = 1.85
$ws.Range("W161").Value = 0.444
$ws.Range("X161").Value = -1
$ws.Range("Y161").Value = -1
$ws.Range("Z161").Value = -0.5
$ws.Range("AA161").Value = 0.4125
$ws.Range("AB161").Value = -1
$ws.Range("AC161").Value = 0.8500000000000001

# Row 162
$ws.Range("A162").Value = 160
$ws.Range("B162").Value = 7032914
$ws.Range("C162").Value = 'Serbia Super Liga'
$ws.Range("D162").Value = 'Serbia Super Liga'
$ws.Range("E162").Value = 45150.57986111111
$ws.Range("F162").Value = 'FK Vozdovac'
$ws.Range("G162").Value = 'FK Radnik Surdulica'
$ws.Range("H162").Value = 1
$ws.Range("I162").Value = 1
$ws.Range("J162").Value = 'D'
$ws.Range("K162").Value = 2.2
$ws.Range("L162").Value = 3.1
$ws.Range("M162").Value = 3.2
$ws.Range("N162").Value = 2.05
$ws.Range("O162").Value = 3.1
$ws.Range("P162").Value = 3.5
$ws.Range("Q162").Value = -0.25
$ws.Range("R162").Value = 1.75
$ws.Range("S162").Value = 2.05
$ws.Range("T162").Value = 2
$ws.Range("U162").Value = 1.775
$ws.Range("V162").Value = 2.025
$ws.Range("W162").Value = -1
$ws.Range("X162").Value = 2.1
$ws.Range("Y162").Value = -1
$ws.Range("Z162").Value = -0.5
$ws.Range("AA162").Value = 0.5249999999999999
$ws.Range("AB162").Value = 0
$ws.Range("AC162").Value = 0

# Row 185
$ws.Range("A185").Value = 183
$ws.Range("B185").Value = 6979440
$ws.Range("C185").Value = 'Serbia Super Liga'
$ws.Range("D185").Value = 'Serbia Super Liga'
$ws.Range("E185").Value = 45171.625
$ws.Range("F185").Value = 'Javor Ivanjica'
$ws.Range("G185").Value = 'Radnicki Nis'
$ws.Range("H185").Value = 1
$ws.Range("I185").Value = 0
$ws.Range("J185").Value = 'H'
$ws.Range("K185").Value = 2.3
$ws.Range("L185").Value = 3.2
$ws.Range("M185").Value = 2.875
$ws.Range("N185").Value = 2.5
$ws.Range("O185").Value = 3.25
$ws.Range("P185").Value = 2.6
$ws.Range("Q185").Value = 0
$ws.Range("R185").Value = 1.85
$ws.Range("S185").Value = 1.95
$ws.Range("T185").Value = 2.25
$ws.Range("U185").Value = 1.9
$ws.Range("V185").Value = 1.9
$ws.Range("W185").Value = 1.5
$ws.Range("X185").Value = -1
$ws.Range("Y185").Value = -1
$ws.Range("Z185").Value = 0.8500000000000001
$ws.Range("AA185").Value = -1
$ws.Range("AB185").Value = -1
$ws.Range("AC185").Value = 0.8999999999999999

# Row 186
$ws.Range("A186").Value = 184
$ws.Range("B186").Value = 6978740
$ws.Range("C186").Value = 'Serbia Super Liga'
$ws.Range("D186").Value = 'Serbia Super Liga'
$ws.Range("E186").Value = 45171.625
$ws.Range("F186").Value = 'Red Star Belgrade'
$ws.Range("G186").Value = 'FK Novi Pazar'
$ws.Range("H186").Value = 2
$ws.Range("I186").Value = 1
$ws.Range("J186").Value = 'H'
$ws.Range("K186").Value = 1.062
$ws.Range("L186").Value = 13
$ws.Range("M186").Value = 23
$ws.Range("N186").Value = 1.025
$ws.Range("O186").Value = 19
$ws.Range("P186").Value = 41
$ws.Range("Q186").Value = -3.75
$ws.Range("R186").Value = 1.825
$ws.Range("S186").Value = 1.975
$ws.Range("T186").Value = 4.5
$ws.Range("U186").Value = 1.975
$ws.Range("V186").Value = 1.825
$ws.Range("W186").Value = 0.02499999999999991
$ws.Range("X186").Value = -1
$ws.Range("Y186").Value = -1
$ws.Range("Z186").Value = -1
$ws.Range("AA186").Value = 0.9750000000000001
$ws.Range("AB186").Value = -1
$ws.Range("AC186").Value = 0.825

# Row 304
$ws.Range("A304").Value = 302
$ws.Range("B304").Value = 6979546
$ws.Range("C304").Value = 'Serbia Super Liga'
$ws.Range("D304").Value = 'Serbia Super Liga'
$ws.Range("E304").Value = 45339.51041666666
$ws.Range("F304").Value = 'Partizan Belgrade'
$ws.Range("G304").Value = 'IMT Novi Belgrade'
$ws.Range("K304").Value = 1.181
$ws.Range("L304").Value = 5.75
$ws.Range("M304").Value = 11
$ws.Range("N304").Value = 1.2
$ws.Range("O304").Value = 5.75
$ws.Range("P304").Value = 9
$ws.Range("Q304").Value = -1.75
$ws.Range("R304").Value = 1.825
$ws.Range("S304").Value = 1.975
$ws.Range("T304").Value = 3.25
$ws.Range("U304").Value = 1.975
$ws.Range("V304").Value = 1.825
$ws.Range("W304").Value = 0
$ws.Range("X304").Value = 0
$ws.Range("Y304").Value = 0
$ws.Range("Z304").Value = 0
$ws.Range("AA304").Value = 0

# Row 305
$ws.Range("A305").Value = 303
$ws.Range("B305").Value = 6979548
$ws.Range("C305").Value = 'Serbia Super Liga'
$ws.Range("D305").Value = 'Serbia Super Liga'
$ws.Range("E305").Value = 45339.60416666666
$ws.Range("F305").Value = 'FK Cukaricki'
$ws.Range("G305").Value = 'FK Zeleznicar Pancevo'
$ws.Range("K305").Value = 1.4
$ws.Range("L305").Value = 4.333
$ws.Range("M305").Value = 6
$ws.Range("N305").Value = 1.4
$ws.Range("O305").Value = 4.333
$ws.Range("P305").Value = 6
$ws.Range("Q305").Value = -1.25
$ws.Range("R305").Value = 1.95
$ws.Range("S305").Value = 1.85
$ws.Range("T305").Value = 2.75
$ws.Range("U305").Value = 1.9
$ws.Range("V305").Value = 1.9
$ws.Range("W305").Value = 0
$ws.Range("X305").Value = 0
$ws.Range("Y305").Value = 0
$ws.Range("Z305").Value = 0
$ws.Range("AA305").Value = 0

# Row 306
$ws.Range("A306").Value = 304
$ws.Range("B306").Value = 6979544
$ws.Range("C306").Value = 'Serbia Super Liga'
$ws.Range("D306").Value = 'Serbia Super Liga'
$ws.Range("E306").Value = 45340.375
$ws.Range("F306").Value = 'FK Radnicki 1923'
$ws.Range("G306").Value = 'Vojvodina'
$ws.Range("K306").Value = 2.8
$ws.Range("L306").Value = 3.25
$ws.Range("M306").Value = 2.2
$ws.Range("N306").Value = 2.9
$ws.Range("O306").Value = 3.25
$ws.Range("P306").Value = 2.1
$ws.Range("Q306").Value = 0.25
$ws.Range("R306").Value = 1.9
$ws.Range("S306").Value = 1.9
$ws.Range("T306").Value = 2.5
$ws.Range("U306").Value = 1.975
$ws.Range("V306").Value = 1.825
$ws.Range("W306").Value = 0
$ws.Range("X306").Value = 0
$ws.Range("Y306").Value = 0
$ws.Range("Z306").Value = 0
$ws.Range("AA306").Value = 0

# Row 307
$ws.Range("A307").Value = 305
$ws.Range("B307").Value = 6979550
$ws.Range("C307").Value = 'Serbia Super Liga'
$ws.Range("D307").Value = 'Serbia Super Liga'
$ws.Range("E307").Value = 45340.45833333334
$ws.Range("F307").Value = 'Spartak Subotica'
$ws.Range("G307").Value = 'Mladost Lucani'
$ws.Range("K307").Value = 2.2
$ws.Range("L307").Value = 3.25
$ws.Range("M307").Value = 2.8
$ws.Range("N307").Value = 2.05
$ws.Range("O307").Value = 3.2
$ws.Range("P307").Value = 3.1
$ws.Range("Q307").Value = -0.25
$ws.Range("R307").Value = 1.8
$ws.Range("S307").Value = 2
$ws.Range("T307").Value = 2.25
$ws.Range("U307").Value = 1.8
$ws.Range("V307").Value = 2
$ws.Range("W307").Value = 0
$ws.Range("X307").Value = 0
$ws.Range("Y307").Value = 0
$ws.Range("Z307").Value = 0
$ws.Range("AA307").Value = 0

# Row 308
$ws.Range("A308").Value = 306
$ws.Range("B308").Value = 6979547
$ws.Range("C308").Value = 'Serbia Super Liga'
$ws.Range("D308").Value = 'Serbia Super Liga'
$ws.Range("E308").Value = 45340.54166666666
$ws.Range("F308").Value = 'FK Backa Topola'
$ws.Range("G308").Value = 'FK Radnik Surdulica'
$ws.Range("K308").Value = 1.333
$ws.Range("L308").Value = 4.333
$ws.Range("M308").Value = 7.5
$ws.Range("N308").Value = 1.333
$ws.Range("O308").Value = 4.333
$ws.Range("P308").Value = 7.5
$ws.Range("Q308").Value = -1.25
$ws.Range("R308").Value = 1.8
$ws.Range("S308").Value = 2
$ws.Range("T308").Value = 2.75
$ws.Range("U308").Value = 1.9
$ws.Range("V308").Value = 1.9
$ws.Range("W308").Value = 0
$ws.Range("X308").Value = 0
$ws.Range("Y308").Value = 0
$ws.Range("Z308").Value = 0
$ws.Range("AA308").Value = 0

# Row 309
$ws.Range("A309").Value = 307
$ws.Range("B309").Value = 6979545
$ws.Range("C309").Value = 'Serbia Super Liga'
$ws.Range("D309").Value = 'Serbia Super Liga'
$ws.Range("E309").Value = 45340.54166666666
$ws.Range("F309").Value = 'Radnicki Nis'
$ws.Range("G309").Value = 'Javor Ivanjica'
$ws.Range("K309").Value = 2
$ws.Range("L309").Value = 3.25
$ws.Range("M309").Value = 3.25
$ws.Range("N309").Value = 1.85
$ws.Range("O309").Value = 3.3
$ws.Range("P309").Value = 3.6
$ws.Range("Q309").Value = -0.5
$ws.Range("R309").Value = 1.85
$ws.Range("S309").Value = 1.95
$ws.Range("T309").Value = 2.5
$ws.Range("U309").Value = 1.975
$ws.Range("V309").Value = 1.825
$ws.Range("W309").Value = 0
$ws.Range("X309").Value = 0
$ws.Range("Y309").Value = 0
$ws.Range("Z309").Value = 0
$ws.Range("AA309").Value = 0

# Row 310
$ws.Range("A310").Value = 308
$ws.Range("B310").Value = 6979549
$ws.Range("C310").Value = 'Serbia Super Liga'
$ws.Range("D310").Value = 'Serbia Super Liga'
$ws.Range("E310").Value = 45341.54166666666
$ws.Range("F310").Value = 'FK Vozdovac'
$ws.Range("G310").Value = 'FK Napredak'
$ws.Range("K310").Value = 2.1
$ws.Range("L310").Value = 3.25
$ws.Range("M310").Value = 3
$ws.Range("N310").Value = 2.1
$ws.Range("O310").Value = 3.25
$ws.Range("P310").Value = 3
$ws.Range("Q310").Value = -0.25
$ws.Range("R310").Value = 1.825
$ws.Range("S310").Value = 1.975
$ws.Range("T310").Value = 2.25
$ws.Range("U310").Value = 1.9
$ws.Range("V310").Value = 1.9
$ws.Range("W310").Value = 0
$ws.Range("X310").Value = 0
$ws.Range("Y310").Value = 0
$ws.Range("Z310").Value = 0
$ws.Range("AA310").Value = 0

# Remove the now-obsolete last row (shift consumed by direct overwrite above)
$ws.Rows(311).Delete() | Out-Null